$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting the existing rows 48-127 down to 49-128.
$ws.Rows(48).Insert()

# Populate the newly inserted row 48 with its data.
$ws.Range("A48").Value = 5
$ws.Range("B48").Value = "Macroferia Regional de Talca"
$ws.Range("C48").Value = "Maule"
$ws.Range("D48").Value = 44725
$ws.Range("E48").Value = 7
$ws.Range("F48").Value = 100112001
$ws.Range("G48").Value = "Berenjena"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 300
$ws.Range("K48").Value = 7000
$ws.Range("L48").Value = 7000
$ws.Range("M48").Value = 7000
$ws.Range("N48").Value = "$/caja 50 unidades"
$ws.Range("O48").Value = "Región de Arica y Parinacota"
$ws.Range("P48").Value = 140
$ws.Range("Q48").Value = 50
$ws.Range("R48").Value = "Hortaliza"
